# Apply the two content changes described by the commit:
#   1. Slide 16's table switches to a different (built-in) table style.
#   2. The deck's theme palette (the "Integral" theme used by the slide
#      master) is replaced with the standard default "Office Theme"
#      palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{AD1DAF26-E188-4DF0-A12C-654BE259E80B}")
    }
}

# --- 2. Theme color palette -------------------------------------------------
function Set-ThemeRGB($colorObj, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorObj.RGB = $r + ($g * 256) + ($b * 65536)
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

# Office theme default color scheme (replaces the "Integral" colors):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeRGB $scheme.Colors(1)  "000000"
Set-ThemeRGB $scheme.Colors(2)  "FFFFFF"
Set-ThemeRGB $scheme.Colors(3)  "44546A"
Set-ThemeRGB $scheme.Colors(4)  "E7E6E6"
Set-ThemeRGB $scheme.Colors(5)  "5B9BD5"
Set-ThemeRGB $scheme.Colors(6)  "ED7D31"
Set-ThemeRGB $scheme.Colors(7)  "A5A5A5"
Set-ThemeRGB $scheme.Colors(8)  "FFC000"
Set-ThemeRGB $scheme.Colors(9)  "4472C4"
Set-ThemeRGB $scheme.Colors(10) "70AD47"
Set-ThemeRGB $scheme.Colors(11) "0563C1"
Set-ThemeRGB $scheme.Colors(12) "954F72"
